$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column K holds the tenant_id / tenant_id_lbl template cells (comment header
# row 1 + model row 2). Import/export should ignore tenant_id, so drop the
# whole column; the following column (update_time_lbl) shifts left into K.
$ws.Range("K1:K2").EntireColumn.Delete()
